$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 400
$ws.Range("I18").Value = 400
$ws.Range("K18").Value = 400
$ws.Range("M18").Value = -116
$ws.Range("H19").Value = 650
$ws.Range("I19").Value = 769.4
$ws.Range("J19").Value = 470.9
$ws.Range("K19").Value = 769.4
$ws.Range("L19").Value = 470.9
$ws.Range("M19").Value = -594.4
$ws.Range("N19").Value = -820.9
$ws.Range("H41").Value = 559.30554
$ws.Range("I41").Value = 719.86957
$ws.Range("J41").Value = 275.23077
$ws.Range("K41").Value = 719.86957
$ws.Range("L41").Value = 275.23077
$ws.Range("M41").Value = -279.86957
$ws.Range("N41").Value = -1155.23077
$ws.Range("H70").Value = 1180
$ws.Range("I70").Value = 1800
$ws.Range("J70").Value = 1025
$ws.Range("K70").Value = 5400
$ws.Range("L70").Value = 3075
$ws.Range("M70").Value = -5130
$ws.Range("N70").Value = -3615
$ws.Range("H73").Value = 1180
$ws.Range("I73").Value = 1800
$ws.Range("J73").Value = 1025
$ws.Range("K73").Value = 5400
$ws.Range("L73").Value = 3075
$ws.Range("M73").Value = -4464
$ws.Range("N73").Value = -4947
$ws.Range("H74").Value = 4133.353
$ws.Range("J74").Value = 4125.3335
$ws.Range("L74").Value = 4125.3335
$ws.Range("N74").Value = -5997.3335
$ws.Range("H76").Value = 4155.273
$ws.Range("I76").Value = 3421.2
$ws.Range("J76").Value = 4767
$ws.Range("K76").Value = 3421.2
$ws.Range("L76").Value = 4767
$ws.Range("M76").Value = -3106.2
$ws.Range("N76").Value = -5397
$ws.Range("H77").Value = 4133.353
$ws.Range("J77").Value = 4125.3335
$ws.Range("L77").Value = 20626.6675
$ws.Range("N77").Value = -29986.6675
$ws.Range("H79").Value = 4155.273
$ws.Range("I79").Value = 3421.2
$ws.Range("J79").Value = 4767
$ws.Range("K79").Value = 3421.2
$ws.Range("L79").Value = 4767
$ws.Range("M79").Value = -2329.2
$ws.Range("N79").Value = -6951
$ws.Range("H100").Value = 2342.9048
$ws.Range("I100").Value = 1625
$ws.Range("K100").Value = 1625
$ws.Range("M100").Value = -1084
$ws.Range("H111").Value = 2082
$ws.Range("I111").Value = 2324
$ws.Range("J111").Value = 1537.5
$ws.Range("K111").Value = 6972
$ws.Range("L111").Value = 4612.5
$ws.Range("M111").Value = -3905
$ws.Range("N111").Value = -10746.5
$ws.Range("H116").Value = 3336266.8
$ws.Range("I116").Value = 3475069.5
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3475069.5
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -3471627.5
$ws.Range("N116").Value = -11884
$ws.Range("H138").Value = 4062.886
$ws.Range("I138").Value = 2130.2
$ws.Range("J138").Value = 4957.648
$ws.Range("K138").Value = 6390.599999999999
$ws.Range("L138").Value = 14872.944
$ws.Range("M138").Value = -1250.599999999999
$ws.Range("N138").Value = -25152.944

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 195
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H43").Value = 6718.375
$ws.Range("J43").Value = 6892.4287
$ws.Range("L43").Value = 6892.4287
$ws.Range("N43").Value = -7518.4287
$ws.Range("H63").Value = 2775
$ws.Range("I63").Value = 2330
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2330
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1644
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 2775
$ws.Range("I66").Value = 2330
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11650
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -8218
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 387943.56
$ws.Range("I74").Value = 3369.0454
$ws.Range("J74").Value = 2503103.5
$ws.Range("K74").Value = 3369.0454
$ws.Range("L74").Value = 2503103.5
$ws.Range("M74").Value = -2495.0454
$ws.Range("N74").Value = -2504851.5
$ws.Range("H77").Value = 387943.56
$ws.Range("I77").Value = 3369.0454
$ws.Range("J77").Value = 2503103.5
$ws.Range("K77").Value = 16845.227
$ws.Range("L77").Value = 12515517.5
$ws.Range("M77").Value = -12477.227
$ws.Range("N77").Value = -12524253.5
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 195
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1957.2162
$ws.Range("I16").Value = 1892.5416
$ws.Range("K16").Value = 1892.5416
$ws.Range("M16").Value = -1605.5416
$ws.Range("H21").Value = 4150
$ws.Range("J21").Value = 6300
$ws.Range("L21").Value = 6300
$ws.Range("N21").Value = -6770
$ws.Range("H31").Value = 2528.6765
$ws.Range("I31").Value = 1481.4762
$ws.Range("J31").Value = 4220.3076
$ws.Range("K31").Value = 1481.4762
$ws.Range("L31").Value = 4220.3076
$ws.Range("M31").Value = -1186.4762
$ws.Range("N31").Value = -4810.3076
$ws.Range("H34").Value = 2528.6765
$ws.Range("I34").Value = 1481.4762
$ws.Range("J34").Value = 4220.3076
$ws.Range("K34").Value = 1481.4762
$ws.Range("L34").Value = 4220.3076
$ws.Range("M34").Value = -1279.4762
$ws.Range("N34").Value = -4624.3076
$ws.Range("H62").Value = 3012.25
$ws.Range("I62").Value = 2260
$ws.Range("J62").Value = 3354.182
$ws.Range("K62").Value = 2260
$ws.Range("L62").Value = 3354.182
$ws.Range("M62").Value = -1636
$ws.Range("N62").Value = -4602.182
$ws.Range("H65").Value = 3012.25
$ws.Range("I65").Value = 2260
$ws.Range("J65").Value = 3354.182
$ws.Range("K65").Value = 11300
$ws.Range("L65").Value = 16770.91
$ws.Range("M65").Value = -8180
$ws.Range("N65").Value = -23010.91
$ws.Range("H113").Value = 1957.2162
$ws.Range("I113").Value = 1892.5416
$ws.Range("K113").Value = 1892.5416
$ws.Range("M113").Value = 277.4584
$ws.Range("H132").Value = 1357.6765
$ws.Range("I132").Value = 1114.5217
$ws.Range("J132").Value = 1866.091
$ws.Range("K132").Value = 3343.5651
$ws.Range("L132").Value = 5598.272999999999
$ws.Range("M132").Value = -813.5650999999998
$ws.Range("N132").Value = -10658.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 694.01514
$ws.Range("I5").Value = 535.45654
$ws.Range("K5").Value = 1606.36962
$ws.Range("M5").Value = -1494.36962
$ws.Range("H135").Value = 694.01514
$ws.Range("I135").Value = 535.45654
$ws.Range("K135").Value = 4819.10886
$ws.Range("M135").Value = -2284.10886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4308.5
$ws.Range("I80").Value = 5500.8335
$ws.Range("J80").Value = 2520
$ws.Range("K80").Value = 5500.8335
$ws.Range("L80").Value = 2520
$ws.Range("M80").Value = -4502.8335
$ws.Range("N80").Value = -4516
$ws.Range("H83").Value = 4308.5
$ws.Range("I83").Value = 5500.8335
$ws.Range("J83").Value = 2520
$ws.Range("K83").Value = 27504.1675
$ws.Range("L83").Value = 12600
$ws.Range("M83").Value = -22512.1675
$ws.Range("N83").Value = -22584
$ws.Range("H97").Value = 2001.6389
$ws.Range("I97").Value = 1433.56
$ws.Range("J97").Value = 3292.7273
$ws.Range("K97").Value = 1433.56
$ws.Range("L97").Value = 3292.7273
$ws.Range("M97").Value = -937.5599999999999
$ws.Range("N97").Value = -4284.7273
$ws.Range("H113").Value = 2444.5833
$ws.Range("I113").Value = 1348
$ws.Range("K113").Value = 1348
$ws.Range("M113").Value = 822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2228.0527
$ws.Range("I61").Value = 1987.3572
$ws.Range("J61").Value = 2902
$ws.Range("K61").Value = 1987.3572
$ws.Range("L61").Value = 2902
$ws.Range("M61").Value = -1785.3572
$ws.Range("N61").Value = -3306
$ws.Range("H113").Value = 2228.0527
$ws.Range("I113").Value = 1987.3572
$ws.Range("J113").Value = 2902
$ws.Range("K113").Value = 1987.3572
$ws.Range("L113").Value = 2902
$ws.Range("M113").Value = 182.6428000000001
$ws.Range("N113").Value = -7242
